$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 280; this shifts the existing rows 280-335
# down to 281-336 (and grows the used range from R335 to R336), matching
# the author's edit of adding one new weekly price record to the table.
$ws.Rows.Item(280).Insert()

# Populate the newly inserted row 280 with the new record. The surrounding
# descriptive columns (market/region/product/variety/quality/unit/origin/
# classification) mirror the row that used to sit at 280 (now at 281);
# only the date, volume, prices and $/Kg differ for this new entry.
$ws.Cells.Item(280, 1).Value2 = 3
$ws.Cells.Item(280, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(280, 3).Value2 = "Coquimbo"
$ws.Cells.Item(280, 4).Value2 = 44637
$ws.Cells.Item(280, 5).Value2 = 5
$ws.Cells.Item(280, 6).Value2 = 100112031
$ws.Cells.Item(280, 7).Value2 = "Poroto verde"
$ws.Cells.Item(280, 8).Value2 = "Magnum"
$ws.Cells.Item(280, 9).Value2 = "Primera"
$ws.Cells.Item(280, 10).Value2 = 115
$ws.Cells.Item(280, 11).Value2 = 26000
$ws.Cells.Item(280, 12).Value2 = 27000
$ws.Cells.Item(280, 13).Value2 = 26504
$ws.Cells.Item(280, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(280, 15).Value2 = "Provincia de Talca"
$ws.Cells.Item(280, 16).Value2 = 1060
$ws.Cells.Item(280, 17).Value2 = 25
$ws.Cells.Item(280, 18).Value2 = "Hortaliza"
